{"js": "const replacements = [\n  [\"2025-03-22 Saturday\", \"2025-03-23 Sunday\"],\n  [\"624\u00d79=5616\", \"466\u00d78=3728\"],\n  [\"523\u00d77=3661\", \"609\u00d73=1827\"],\n  [\"550\u00d77=3850\", \"411\u00d75=2055\"],\n  [\"466\u00d77=3262\", \"550\u00d78=4400\"],\n  [\"452\u00d72=904\", \"959\u00d79=8631\"],\n  [\"238\u00d72=476\", \"609\u00d76=3654\"],\n  [\"732\u00d78=5856\", \"945\u00d73=2835\"],\n  [\"486\u00d75=2430\", \"907\u00d74=3628\"],\n  [\"260\u00d75=1300\", \"143\u00d78=1144\"],\n  [\"278\u00d79=2502\", \"838\u00d73=2514\"],\n  [\"744\u00d74=2976\", \"563\u00d75=2815\"],\n  [\"822\u00d76=4932\", \"195\u00d76=1170\"],\n  [\"273\u00d77=1911\", \"943\u00d77=6601\"],\n  [\"539\u00d78=4312\", \"383\u00d73=1149\"],\n  [\"796\u00d78=6368\", \"524\u00d75=2620\"],\n  [\"843\u00d72=1686\", \"558\u00d77=3906\"],\n  [\"357\u00d75=1785\", \"396\u00d74=1584\"],\n  [\"852\u00d76=5112\", \"396\u00d75=1980\"],\n  [\"464\u00d75=2320\", \"273\u00d78=2184\"],\n  [\"727\u00d74=2908\", \"660\u00d72=1320\"],\n  [\"698\u00d76=4188\", \"241\u00d75=1205\"],\n  [\"526\u00d77=3682\", \"424\u00d72=848\"],\n  [\"725\u00d75=3625\", \"136\u00d76=816\"],\n  [\"993\u00d76=5958\", \"970\u00d76=5820\"],\n  [\"748\u00d72=1496\", \"483\u00d73=1449\"],\n];\n\nconst body = context.document.body;\nlet totalReplaced = 0;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, {matchCase: true, matchWholeWord: false});\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n    totalReplaced++;\n  }\n  await context.sync();\n}\n\nreturn \"replaced \" + totalReplaced + \" of \" + replacements.length + \" entries\";\n", "ps1": "$d = $word.ActiveDocument\n$replacements = @(\n    @('2025-03-22 Saturday', '2025-03-23 Sunday'),\n    @('624\u00d79=5616', '466\u00d78=3728'),\n    @('523\u00d77=3661', '609\u00d73=1827'),\n    @('550\u00d77=3850', '411\u00d75=2055'),\n    @('466\u00d77=3262', '550\u00d78=4400'),\n    @('452\u00d72=904', '959\u00d79=8631'),\n    @('238\u00d72=476', '609\u00d76=3654'),\n    @('732\u00d78=5856', '945\u00d73=2835'),\n    @('486\u00d75=2430', '907\u00d74=3628'),\n    @('260\u00d75=1300', '143\u00d78=1144'),\n    @('278\u00d79=2502', '838\u00d73=2514'),\n    @('744\u00d74=2976', '563\u00d75=2815'),\n    @('822\u00d76=4932', '195\u00d76=1170'),\n    @('273\u00d77=1911', '943\u00d77=6601'),\n    @('539\u00d78=4312', '383\u00d73=1149'),\n    @('796\u00d78=6368', '524\u00d75=2620'),\n    @('843\u00d72=1686', '558\u00d77=3906'),\n    @('357\u00d75=1785', '396\u00d74=1584'),\n    @('852\u00d76=5112', '396\u00d75=1980'),\n    @('464\u00d75=2320', '273\u00d78=2184'),\n    @('727\u00d74=2908', '660\u00d72=1320'),\n    @('698\u00d76=4188', '241\u00d75=1205'),\n    @('526\u00d77=3682', '424\u00d72=848'),\n    @('725\u00d75=3625', '136\u00d76=816'),\n    @('993\u00d76=5958', '970\u00d76=5820'),\n    @('748\u00d72=1496', '483\u00d73=1449'),\n)\n\n$count = 0\nforeach ($pair in $replacements) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $found = $rng.Find.Execute($pair[0], $true, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n    if ($found) { $count++ }\n}\nWrite-Output \"replaced $count of $($replacements.Count) entries\"\n"}
